# Daily attendance processing - 2026-01-08 21:34:53
# Swap the order of names in the "Recorded By" column (G) from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
